# Fiscal Sponsorship Agreement - targeted text edits
#
# NOTE: Find.Execute(...,Replace:=...) silently "smart-quotifies" straight
# apostrophes into curly ones in this runtime even with AutoFormat options
# off, so every replacement below locates the target with Find.Execute()
# (no replacement arg) and then assigns Range.Text directly, which leaves
# straight apostrophes untouched.

$d = $word.ActiveDocument

# --- Change 1: Section 1 - remove "programmatic, " and collapse the double spaces ---
$r1 = $d.Content
$r1.Find.Text = "1. The Fiscal Sponsor hereby agrees to sponsor the Project and to assume administrative,  programmatic, financial, and legal responsibility for purposes of the requirements of  funding organizations. The Sponsored Organization agrees to implement and operate the  Project, in accordance with the terms of this agreement and with any requirements imposed  by funding organizations.  "
$r1.Find.MatchCase = $true
$r1.Find.MatchWildcards = $false
$found1 = $r1.Find.Execute()
if ($found1) {
    $r1.Text = "1. The Fiscal Sponsor hereby agrees to sponsor the Project and to assume administrative, financial, and legal responsibility for purposes of the requirements of funding organizations. The Sponsored Organization agrees to implement and operate the Project, in accordance with the terms of this agreement and with any requirements imposed by funding organizations. "
}
Write-Output ("change1 found=" + $found1)

# --- Change 2: Section 14 - "program design" -> "programmatic design"; rewrite last sentence ---
$r2 = $d.Content
$r2.Find.Text = "14. The Sponsored Organization retains decision-making authority over Project strategy, program design, partnerships with Ukrainian municipalities, hiring of contractors or staff, and day-to-day operations, subject to the Fiscal Sponsor's reserved variance power under Section 15, provided such decisions: (a) comply with Section 501(c)(3) requirements, (b) do not jeopardize the Fiscal Sponsor's tax-exempt status, and (c) are consistent with approved grant terms where applicable. The Fiscal Sponsor's role is to provide administrative and fiduciary support, and the parties recognize that programmatic decisions rest with the Project team."
$r2.Find.MatchCase = $true
$r2.Find.MatchWildcards = $false
$found2 = $r2.Find.Execute()
if ($found2) {
    $r2.Text = "14. The Sponsored Organization retains decision-making authority over Project strategy, programmatic design, partnerships with Ukrainian municipalities, hiring of contractors or staff, and day-to-day operations, subject to the Fiscal Sponsor's reserved variance power under Section 15, provided such decisions: (a) comply with Section 501(c)(3) requirements, (b) do not jeopardize the Fiscal Sponsor's tax-exempt status, and (c) are consistent with approved grant terms where applicable. The Fiscal Sponsor provides administrative and fiduciary oversight, and the parties intend for the Project team to manage day-to-day operations within the framework of that oversight."
}
Write-Output ("change2 found=" + $found2)

# --- Change 3: variance-power paragraph under Section 15 - drop "programmatic " and delete last sentence ---
$r3 = $d.Content
$r3.Find.Text = "The Fiscal Sponsor's exercise of variance power under this section is limited to ensuring compliance with applicable law, IRS requirements, U.S. sanctions laws and regulations administered by the Office of Foreign Assets Control (OFAC), and the Fiscal Sponsor's tax-exempt purposes. Subject to this reserved authority, the Sponsored Organization retains full programmatic and operational autonomy as described in Section 14. The Fiscal Sponsor shall not exercise variance power to override the Sponsored Organization's programmatic decisions except where necessary to maintain compliance with legal, regulatory, or sanctions requirements."
$r3.Find.MatchCase = $true
$r3.Find.MatchWildcards = $false
$found3 = $r3.Find.Execute()
if ($found3) {
    $r3.Text = "The Fiscal Sponsor's exercise of variance power under this section is limited to ensuring compliance with applicable law, IRS requirements, U.S. sanctions laws and regulations administered by the Office of Foreign Assets Control (OFAC), and the Fiscal Sponsor's tax-exempt purposes. Subject to this reserved authority, the Sponsored Organization retains full and operational autonomy as described in Section 14."
}
Write-Output ("change3 found=" + $found3)

# --- Change 4: insert a new Section 24 (Force Majeure); old Section 24 becomes Section 25 ---
$r4 = $d.Content
$r4.Find.Text = "24. This Agreement shall be governed by and construed in accordance with the laws of the Commonwealth of Virginia, without regard to its conflict of laws provisions."
$r4.Find.MatchCase = $true
$r4.Find.MatchWildcards = $false
$found4 = $r4.Find.Execute()
Write-Output ("change4 found=" + $found4)
if ($found4) {
    # Remember the old text, then overwrite the paragraph with the Force Majeure
    # clause (keeping the "24." numbering), and push the original sentence into
    # a freshly inserted paragraph right after it, renumbered to "25.".
    $oldSection24Text = "25. This Agreement shall be governed by and construed in accordance with the laws of the Commonwealth of Virginia, without regard to its conflict of laws provisions."
    $r4.Text = "24. Force Majeure. Neither party shall be liable for any failure or delay in performing its obligations under this Agreement to the extent such failure or delay results from circumstances beyond its reasonable control, including but not limited to: armed conflict or military operations affecting Ukraine, changes in U.S. or international sanctions, disruption of international banking or wire transfer channels, natural disasters, pandemics, or actions by governmental authorities. The affected party shall notify the other party promptly and use reasonable efforts to mitigate the impact. Obligations shall resume once the force majeure condition ceases. If a force majeure event prevents performance for more than ninety (90) consecutive days, either party may terminate this Agreement upon thirty (30) days written notice, subject to the transition provisions of Section 12."

    $r4.InsertParagraphAfter()
    $newParaStart = $r4.End + 1
    $newPara = $d.Range($newParaStart, $newParaStart)
    $newPara.Text = $oldSection24Text
}

# --- Change 5: renumber the old Section 25 ("entire understanding...") to Section 26 ---
$r5 = $d.Content
$r5.Find.Text = "25. This Agreement constitutes the entire understanding between the parties and supersedes all prior agreements relating to its subject matter. This Agreement may be amended only by written instrument signed by both parties."
$r5.Find.MatchCase = $true
$r5.Find.MatchWildcards = $false
$found5 = $r5.Find.Execute()
if ($found5) {
    $r5.Text = "26. This Agreement constitutes the entire understanding between the parties and supersedes all prior agreements relating to its subject matter. This Agreement may be amended only by written instrument signed by both parties."
}
Write-Output ("change5 found=" + $found5)

Write-Output "done"
